$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 604.86
$ws.Range("C3").Value = 597.38
$ws.Range("C4").Value = 585.78
$ws.Range("C5").Value = 586.49
$ws.Range("C6").Value = 586.49
